# Auto-generated edit script: updates market-price derived columns (H-N)
# on the Zeromus_Profits workbook per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 820
$ws.Range("I28").Value = 852.5
$ws.Range("J28").Value = 690
$ws.Range("K28").Value = 852.5
$ws.Range("L28").Value = 690
$ws.Range("M28").Value = -367.5
$ws.Range("N28").Value = -1660
$ws.Range("H33").Value = 742.8857400000001
$ws.Range("I33").Value = 808.6774
$ws.Range("K33").Value = 808.6774
$ws.Range("M33").Value = -579.6774
$ws.Range("H40").Value = 3267.9167
$ws.Range("I40").Value = 3767.0715
$ws.Range("J40").Value = 2569.1
$ws.Range("K40").Value = 3767.0715
$ws.Range("L40").Value = 2569.1
$ws.Range("M40").Value = -3592.0715
$ws.Range("N40").Value = -2919.1
$ws.Range("H64").Value = 3014.7036
$ws.Range("J64").Value = 3037.3157
$ws.Range("L64").Value = 3037.3157
$ws.Range("N64").Value = -3533.3157
$ws.Range("H67").Value = 3014.7036
$ws.Range("J67").Value = 3037.3157
$ws.Range("L67").Value = 3037.3157
$ws.Range("N67").Value = -4753.3157
$ws.Range("H113").Value = 2969
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3043.5386
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3043.5386
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -9551.5386
$ws.Range("H116").Value = 2139645.2
$ws.Range("I116").Value = 2453819
$ws.Range("K116").Value = 2453819
$ws.Range("M116").Value = -2450377
$ws.Range("H130").Value = 39779.875
$ws.Range("J130").Value = 39779.875
$ws.Range("L130").Value = 39779.875
$ws.Range("N130").Value = -49819.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1192.2
$ws.Range("I2").Value = 987
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 987
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -874
$ws.Range("N2").Value = -1726
$ws.Range("H21").Value = 1608.5
$ws.Range("I21").Value = 1200
$ws.Range("J21").Value = 2017
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 2017
$ws.Range("M21").Value = -826
$ws.Range("N21").Value = -2765
$ws.Range("H32").Value = 21017.658
$ws.Range("I32").Value = 6797.1665
$ws.Range("K32").Value = 6797.1665
$ws.Range("M32").Value = -6510.1665
$ws.Range("H45").Value = 2434.9333
$ws.Range("I45").Value = 2758.2222
$ws.Range("J45").Value = 1950
$ws.Range("K45").Value = 2758.2222
$ws.Range("L45").Value = 1950
$ws.Range("M45").Value = -2381.2222
$ws.Range("N45").Value = -2704
$ws.Range("H61").Value = 2040.2609
$ws.Range("I61").Value = 1684.6666
$ws.Range("J61").Value = 2428.182
$ws.Range("K61").Value = 1684.6666
$ws.Range("L61").Value = 2428.182
$ws.Range("M61").Value = -1472.6666
$ws.Range("N61").Value = -2852.182
$ws.Range("H74").Value = 4732.3335
$ws.Range("I74").Value = 7020.4287
$ws.Range("J74").Value = 3276.2727
$ws.Range("K74").Value = 7020.4287
$ws.Range("L74").Value = 3276.2727
$ws.Range("M74").Value = -6146.4287
$ws.Range("N74").Value = -5024.2727
$ws.Range("H77").Value = 4732.3335
$ws.Range("I77").Value = 7020.4287
$ws.Range("J77").Value = 3276.2727
$ws.Range("K77").Value = 35102.14350000001
$ws.Range("L77").Value = 16381.3635
$ws.Range("M77").Value = -30734.14350000001
$ws.Range("N77").Value = -25117.3635
$ws.Range("H110").Value = 1871.32
$ws.Range("I110").Value = 1464.6
$ws.Range("J110").Value = 3498.2
$ws.Range("K110").Value = 1464.6
$ws.Range("L110").Value = 3498.2
$ws.Range("M110").Value = 580.4000000000001
$ws.Range("N110").Value = -7588.2
$ws.Range("H116").Value = 1192.2
$ws.Range("I116").Value = 987
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 987
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1307
$ws.Range("N116").Value = -6088
$ws.Range("H133").Value = 26755.715
$ws.Range("J133").Value = 26755.715
$ws.Range("L133").Value = 26755.715
$ws.Range("N133").Value = -31815.715
$ws.Range("H136").Value = 2040.2609
$ws.Range("I136").Value = 1684.6666
$ws.Range("J136").Value = 2428.182
$ws.Range("K136").Value = 5053.9998
$ws.Range("L136").Value = 7284.545999999999
$ws.Range("M136").Value = -2503.9998
$ws.Range("N136").Value = -12384.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1192.2
$ws.Range("I3").Value = 987
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 987
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -873
$ws.Range("N3").Value = -1728

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 71429800
$ws.Range("I16").Value = 166667500
$ws.Range("J16").Value = 1522.5
$ws.Range("K16").Value = 166667500
$ws.Range("L16").Value = 1522.5
$ws.Range("M16").Value = -166667213
$ws.Range("N16").Value = -2096.5
$ws.Range("H31").Value = 4764084.5
$ws.Range("I31").Value = 10001166
$ws.Range("J31").Value = 3101.5454
$ws.Range("K31").Value = 10001166
$ws.Range("L31").Value = 3101.5454
$ws.Range("M31").Value = -10000871
$ws.Range("N31").Value = -3691.5454
$ws.Range("H34").Value = 4764084.5
$ws.Range("I34").Value = 10001166
$ws.Range("J34").Value = 3101.5454
$ws.Range("K34").Value = 10001166
$ws.Range("L34").Value = 3101.5454
$ws.Range("M34").Value = -10000964
$ws.Range("N34").Value = -3505.5454
$ws.Range("H107").Value = 19231604
$ws.Range("I107").Value = 27778456
$ws.Range("K107").Value = 27778456
$ws.Range("M107").Value = -27776536
$ws.Range("H113").Value = 71429800
$ws.Range("I113").Value = 166667500
$ws.Range("J113").Value = 1522.5
$ws.Range("K113").Value = 166667500
$ws.Range("L113").Value = 1522.5
$ws.Range("M113").Value = -166665330
$ws.Range("N113").Value = -5862.5
$ws.Range("H132").Value = 2433.7896
$ws.Range("I132").Value = 1952.75
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 5858.25
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -3328.25
$ws.Range("N132").Value = -20058.0005
$ws.Range("H134").Value = 4731.12
$ws.Range("I134").Value = 4810.5293
$ws.Range("J134").Value = 4562.375
$ws.Range("K134").Value = 14431.5879
$ws.Range("L134").Value = 13687.125
$ws.Range("M134").Value = -11896.5879
$ws.Range("N134").Value = -18757.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 328.1875
$ws.Range("I8").Value = 328.1875
$ws.Range("K8").Value = 984.5625
$ws.Range("M8").Value = -845.5625
$ws.Range("H40").Value = 7700.077
$ws.Range("J40").Value = 12471.375
$ws.Range("L40").Value = 49885.5
$ws.Range("N40").Value = -50023.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 2853.5
$ws.Range("I23").Value = 400
$ws.Range("J23").Value = 3076.5454
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 3076.5454
$ws.Range("M23").Value = -177
$ws.Range("N23").Value = -3522.5454
$ws.Range("H24").Value = 20802004
$ws.Range("I24").Value = 26000004
$ws.Range("J24").Value = 10007
$ws.Range("K24").Value = 26000004
$ws.Range("L24").Value = 10007
$ws.Range("M24").Value = -25999831
$ws.Range("N24").Value = -10353
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 9914.286
$ws.Range("J40").Value = 9914.286
$ws.Range("L40").Value = 9914.286
$ws.Range("N40").Value = -10216.286
$ws.Range("H44").Value = 9950
$ws.Range("J44").Value = 9950
$ws.Range("L44").Value = 9950
$ws.Range("N44").Value = -11142
$ws.Range("H46").Value = 4100
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H57").Value = 20061
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 20061
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 20061
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -21701
$ws.Range("H58").Value = 17800
$ws.Range("J58").Value = 17800
$ws.Range("L58").Value = 17800
$ws.Range("N58").Value = -18354
$ws.Range("H113").Value = 1475
$ws.Range("I113").Value = 1112.5
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1112.5
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 1057.5
$ws.Range("N113").Value = -6540
$ws.Range("H132").Value = 2328.7585
$ws.Range("I132").Value = 1475.6316
$ws.Range("J132").Value = 3949.7
$ws.Range("K132").Value = 4426.8948
$ws.Range("L132").Value = 11849.1
$ws.Range("M132").Value = -1896.8948
$ws.Range("N132").Value = -16909.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 7505
$ws.Range("J5").Value = 7505
$ws.Range("L5").Value = 7505
$ws.Range("N5").Value = -7731

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9996
$ws.Range("J2").Value = 9996
$ws.Range("L2").Value = 9996
$ws.Range("N2").Value = -10220
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H23").Value = 1410
$ws.Range("I23").Value = 1410
$ws.Range("K23").Value = 1410
$ws.Range("M23").Value = -1181
$ws.Range("H25").Value = 9021.6
$ws.Range("J25").Value = 9021.6
$ws.Range("L25").Value = 9021.6
$ws.Range("N25").Value = -9607.6
$ws.Range("H136").Value = 1731.2
$ws.Range("I136").Value = 1164.5714
$ws.Range("J136").Value = 2227
$ws.Range("K136").Value = 3493.7142
$ws.Range("L136").Value = 6681
$ws.Range("M136").Value = -943.7142000000003
$ws.Range("N136").Value = -11781
